$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44342
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 250
$ws.Range("N2").Value = 24000
$ws.Range("O2").Value = 25000
$ws.Range("P2").Value = 24500
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1361

# Row 3
$ws.Range("D3").Value = 44355
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 270
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 21000
$ws.Range("P3").Value = 20500
$ws.Range("R3").Value = "Región Metropolitana"
$ws.Range("S3").Value = 1139

# Row 4
$ws.Range("D4").Value = 44301
$ws.Range("K4").Value = "Hachiya"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 21000
$ws.Range("P4").Value = 20500
$ws.Range("S4").Value = 1139

# Row 6
$ws.Range("D6").Value = 44313
$ws.Range("K6").Value = "Mankaki"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 270
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 21500
$ws.Range("S6").Value = 1194
